# Applies the "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
# edit: refresh the account-statement figures and expand the worker/period
# detail table from 3 rows to 5 rows (3 distinct workers, 4 distinct periods).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# ---------------------------------------------------------------------------
# 1. Header summary figures
# ---------------------------------------------------------------------------
$ws.Range("E11").Value2 = 172814   # VALOR MORA total
$ws.Range("C13").Value2 = 3        # Cant. Trabajadores
$ws.Range("F13").Value2 = 4        # Cant. Periodos

# ---------------------------------------------------------------------------
# 2. Grow the detail table from 3 data rows (16-18) to 5 data rows (16-20).
#    Row 18 currently carries the "last row" (bottom border) style, so two
#    blank rows are inserted above it -- this pushes row 18 down to row 20
#    (keeping its special style), while the new rows 18/19 are given the
#    regular "middle" row style copied from row 17. Everything below (the
#    footer block) shifts down automatically with the insert.
# ---------------------------------------------------------------------------
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(18).Insert()

$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B17:J17").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Fill in the detail rows with the updated worker / period data.
# ---------------------------------------------------------------------------
# Row 16: BLANCA ROSA MOYA NARVAEZ - periodo 2505
$ws.Range("B16").Value2 = "CC"
$ws.Range("C16").Value2 = "45556094"
$ws.Range("D16").Value2 = "BLANCA ROSA MOYA NARVAEZ"
$ws.Range("E16").Value2 = "2505"
$ws.Range("F16").Value2 = 56940
$ws.Range("G16").Value2 = 1423500

# Row 17: BLANCA ROSA MOYA NARVAEZ - periodo 2412
$ws.Range("B17").Value2 = "CC"
$ws.Range("C17").Value2 = "45556094"
$ws.Range("D17").Value2 = "BLANCA ROSA MOYA NARVAEZ"
$ws.Range("E17").Value2 = "2412"
$ws.Range("F17").Value2 = 52000
$ws.Range("G17").Value2 = 1423500

# Row 18: BLANCA ROSA MOYA NARVAEZ - periodo 2410
$ws.Range("B18").Value2 = "CC"
$ws.Range("C18").Value2 = "45556094"
$ws.Range("D18").Value2 = "BLANCA ROSA MOYA NARVAEZ"
$ws.Range("E18").Value2 = "2410"
$ws.Range("F18").Value2 = 3467
$ws.Range("G18").Value2 = 1423500

# Row 19: SANTIAGO MOYA NARVAEZ - periodo 2501
$ws.Range("B19").Value2 = "CC"
$ws.Range("C19").Value2 = "1047401502"
$ws.Range("D19").Value2 = "SANTIAGO MOYA NARVAEZ"
$ws.Range("E19").Value2 = "2501"
$ws.Range("F19").Value2 = 56940
$ws.Range("G19").Value2 = 1300000

# Row 20: EFEMBER GONZALEZ MOYA - periodo 2410 (keeps the "last row" style)
$ws.Range("B20").Value2 = "CC"
$ws.Range("C20").Value2 = "1007138237"
$ws.Range("D20").Value2 = "EFEMBER GONZALEZ MOYA"
$ws.Range("E20").Value2 = "2410"
$ws.Range("F20").Value2 = 3467
$ws.Range("G20").Value2 = 1300000
